$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...Personal Data Ecosystems [2.3.4]..." becomes
#           "...Personal Data Ecosystems and MyData [2.3.4]..."
#           with "MyData" in italics, inserted right before the existing
#           hyperlinked "[" that opens the "2.3.4" citation.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Personal Data Ecosystems [", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Personal Data Ecosystems [' text to update"
}

$matchStart = $rng.Start

# Replace the matched text (which ends in the opening "[" of the citation)
# with just the new lead-in text (no trailing bracket yet).
$leadIn = "Personal Data Ecosystems and"
$rng.Text = $leadIn
$pos = $matchStart + $leadIn.Length

# " " (plain)
$d.Range($pos, $pos).Text = " "
$pos = $pos + 1

# "MyData" (italic)
$d.Range($pos, $pos).Text = "MyData"
$myDataStart = $pos
$myDataEnd = $pos + 6
$pos = $myDataEnd

# " " (plain)
$d.Range($pos, $pos).Text = " "
$pos = $pos + 1

# "[" (plain) - re-adds the bracket that used to directly follow "Ecosystems"
$d.Range($pos, $pos).Text = "["
$pos = $pos + 1

# Apply italic formatting to just the "MyData" run
$d.Range($myDataStart, $myDataEnd).Font.Italic = $true

# ---------------------------------------------------------------------------
# Change 2: Tweak the closing sentence of the same paragraph.
#   "]. Collectively, I now knew where to position my existing and newly
#   discovered understandings..." becomes
#   "]. Collectively through these discoveries, I learned knew where to
#   position my existing and newly-discovered understandings..."
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(". Collectively, I now knew where to position", $true, $false, $false, $false, $false, $true, 1, $false, ". Collectively through these discoveries, I learned knew where to position", 2)
if (-not $found2) {
    throw "Could not find the 'Collectively, I now knew' sentence to update"
}

$rng3 = $d.Content
$found3 = $rng3.Find.Execute("newly discovered understandings", $true, $false, $false, $false, $false, $true, 1, $false, "newly-discovered understandings", 2)
if (-not $found3) {
    throw "Could not find 'newly discovered understandings' text to update"
}
